$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Sigma-X (Triang)"
$ws.Range("C1").Value = "Sigma-Y (Triang)"
$ws.Range("D1").Value = "Tau-XY (Triang)"
$ws.Range("E1").Value = "Sigma-VM (Triang)"

$ws.Columns.Item(2).ColumnWidth = 15.5703125
$ws.Columns.Item(3).ColumnWidth = 15.42578125
$ws.Columns.Item(5).ColumnWidth = 17.42578125
